# Natmi following Dr Hou advice
# Update rows 2-4 (existing Sending/Ligand/Receptor combos) with recalculated
# statistics, and append two new rows (ECs and self-loop sCs target clusters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: sCs / Il12a / Il12rb2 / FAPs ---
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 0.8401833333333332
$ws.Cells.Item(2, 8).Value = 2.52055
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1353843333333333
$ws.Cells.Item(2, 14).Value = 0.406153
$ws.Cells.Item(2, 15).Value = 0.05860584767036442
$ws.Cells.Item(2, 16).Value = 0.05860584767036442
$ws.Cells.Item(2, 17).Value = 0.1137476604611111
$ws.Cells.Item(2, 18).Value = 1.02372894415
$ws.Cells.Item(2, 19).Value = 0.05860584767036442
$ws.Cells.Item(2, 20).Value = 0.05860584767036442

# --- Row 3: sCs / Il12a / Il12rb2 / M1 ---
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 0.8401833333333332
$ws.Cells.Item(3, 8).Value = 2.52055
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.5252536666666666
$ws.Cells.Item(3, 14).Value = 1.575761
$ws.Cells.Item(3, 15).Value = 0.2273744355720654
$ws.Cells.Item(3, 16).Value = 0.2273744355720655
$ws.Cells.Item(3, 17).Value = 0.4413093765055554
$ws.Cells.Item(3, 18).Value = 3.971784388549999
$ws.Cells.Item(3, 19).Value = 0.2273744355720654
$ws.Cells.Item(3, 20).Value = 0.2273744355720655

# --- Row 4: sCs / Il12a / Il12rb2 / M2 ---
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 0.8401833333333332
$ws.Cells.Item(4, 8).Value = 2.52055
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.390367
$ws.Cells.Item(4, 14).Value = 1.171101
$ws.Cells.Item(4, 15).Value = 0.168984020338669
$ws.Cells.Item(4, 16).Value = 0.168984020338669
$ws.Cells.Item(4, 17).Value = 0.3279798472833332
$ws.Cells.Item(4, 18).Value = 2.95181862555
$ws.Cells.Item(4, 19).Value = 0.168984020338669
$ws.Cells.Item(4, 20).Value = 0.168984020338669

# --- Row 5 (new): sCs / Il12a / Il12rb2 / ECs ---
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Il12a"
$ws.Cells.Item(5, 3).Value = "Il12rb2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.8401833333333332
$ws.Cells.Item(5, 8).Value = 2.52055
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.916887
$ws.Cells.Item(5, 14).Value = 2.750661
$ws.Cells.Item(5, 15).Value = 0.3969066326207421
$ws.Cells.Item(5, 16).Value = 0.3969066326207421
$ws.Cells.Item(5, 17).Value = 0.7703531759499999
$ws.Cells.Item(5, 18).Value = 6.933178583549999
$ws.Cells.Item(5, 19).Value = 0.3969066326207421
$ws.Cells.Item(5, 20).Value = 0.3969066326207421

# --- Row 6 (new): sCs / Il12a / Il12rb2 / sCs (self-loop) ---
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Il12a"
$ws.Cells.Item(6, 3).Value = "Il12rb2"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.8401833333333332
$ws.Cells.Item(6, 8).Value = 2.52055
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.3421903333333333
$ws.Cells.Item(6, 14).Value = 1.026571
$ws.Cells.Item(6, 15).Value = 0.148129063798159
$ws.Cells.Item(6, 16).Value = 0.148129063798159
$ws.Cells.Item(6, 17).Value = 0.2875026148944444
$ws.Cells.Item(6, 18).Value = 2.587523534049999
$ws.Cells.Item(6, 19).Value = 0.148129063798159
$ws.Cells.Item(6, 20).Value = 0.148129063798159
